$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Ost/Nord coordinates to rounded integer values
$ws.Range("Q2").Value = 655459
$ws.Range("R2").Value = 7218293

# Clear Starttid (Z2) and Sluttid (AB2) - these were inline string "00:00" cells
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
